$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# Force text format on Price cells so numeric-looking strings (e.g. "8.00", "0.999")
# keep their exact textual representation instead of being converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.192.74"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.450.24"

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.13"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.47"
$ws.Range("E6").Value = "  -1.84%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.449.38"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -3.60%  "

$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("E13").Value = "  -3.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.49"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("E15").Value = "  -3.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.199.39"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.429.98"
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("E19").Value = "  -3.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.15"
$ws.Range("E20").Value = "  -2.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.97"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("E22").Value = "  -2.71%  "

$ws.Range("E23").Value = "  -3.82%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.89"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.39"
$ws.Range("E26").Value = "  +6.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "617.97"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0957"
$ws.Range("E28").Value = "  -6.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.540.56"
$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("E31").Value = "  -4.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"

$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("E35").Value = "  -5.34%  "

$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("E37").Value = "  -5.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "150.26"
$ws.Range("E39").Value = "  +2.48%  "

$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.35"
$ws.Range("E41").Value = "  -2.65%  "

$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.53"
$ws.Range("E43").Value = "  +1.85%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("E45").Value = "  -4.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.22"
$ws.Range("E46").Value = "  -3.63%  "

$ws.Range("E47").Value = "  -3.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.606"
$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("E50").Value = "  +12.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.55"
$ws.Range("E51").Value = "  -7.38%  "

Write-Host "Updated cryptos list values."
